$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.02%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.99%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.492"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.21%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08001"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.84%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.974"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.82%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.577"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-4.67%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.50%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1110"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.70%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1909"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.58%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "9.571"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "13.44%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09980"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.04%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04777"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12.77%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1063"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.34%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001271"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.78%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04078"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.18%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005936"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.07%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.34%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.388"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.03%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.26%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.63%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2584"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.98%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001272"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.46%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004367"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.10%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.52%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003743"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.24%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02597"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.34%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05732"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.85%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007546"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.89%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1401"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.14%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.91%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002014"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.20%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008338"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.84%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007123"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.10%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.11%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005798"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.23%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003529"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "55.34%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.47%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.11%"
